# Update "想去人数" (F column) counts on the sheets "展览" and "全部类型"
# to match the newly generated data output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Map of row -> new F value for sheet "展览"
$exhibitionUpdates = @{
    2  = 67
    3  = 515
    6  = 364
    8  = 2215
    10 = 5443
    11 = 125
    12 = 361
}

# Map of row -> new F value for sheet "全部类型"
$allTypesUpdates = @{
    2  = 67
    4  = 515
    7  = 364
    11 = 2215
    13 = 5443
    14 = 125
    15 = 361
}

$wsExhibition = $wb.Worksheets.Item("展览")
foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibition.Cells.Item($row, 6).Value = $exhibitionUpdates[$row]
}

$wsAllTypes = $wb.Worksheets.Item("全部类型")
foreach ($row in $allTypesUpdates.Keys) {
    $wsAllTypes.Cells.Item($row, 6).Value = $allTypesUpdates[$row]
}
